$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 169, pushing the existing rows 169-172 down to 170-173.
$ws.Rows.Item(169).Insert()

# Populate the newly inserted row 169 with this week's record (matches the
# existing pattern for this Mercado/Categoria block).
$ws.Range("A169").Value = 10
$ws.Range("B169").Value = "Vega Modelo de Temuco"
$ws.Range("C169").Value = "La Araucanía"
$ws.Range("D169").Value = 44509
$ws.Range("D169").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E169").Value = 9
$ws.Range("F169").Value = 100112039
$ws.Range("G169").Value = "Ciboulette"
$ws.Range("H169").Value = "Sin especificar"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 20
$ws.Range("K169").Value = 5000
$ws.Range("L169").Value = 5000
$ws.Range("M169").Value = 5000
$ws.Range("N169").Value = "$/docena de atados"
$ws.Range("O169").Value = "Provincia de Cautín"
$ws.Range("P169").Value = 1667
$ws.Range("Q169").Value = 3
$ws.Range("R169").Value = "Hortaliza"
